$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column B for rows 2-10 (stock tickers removed from "Buying Opportunity" column)
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()

# Clear D2 and F2 (long buildup / FII entering no longer populated)
$ws.Range("D2").ClearContents()
$ws.Range("F2").ClearContents()

# Update column C (support Zone) values for existing rows 2-28
$ws.Range("C2").Value = "NSE:ALKYLAMINE"
$ws.Range("C3").Value = "NSE:AMBIKCO"
$ws.Range("C4").Value = "NSE:ASIANHOTNR"
$ws.Range("C5").Value = "NSE:BARBEQUE"
$ws.Range("C6").Value = "NSE:BIRLACORPN"
$ws.Range("C7").Value = "NSE:CHEMFAB"
$ws.Range("C8").Value = "NSE:CINELINE"
$ws.Range("C9").Value = "NSE:CYBERMEDIA"
$ws.Range("C10").Value = "NSE:EXCELINDUS"
$ws.Range("C11").Value = "NSE:FOSECOIND"
$ws.Range("C12").Value = "NSE:GLAXO"
$ws.Range("C13").Value = "NSE:GMMPFAUDLR"
$ws.Range("C14").Value = "NSE:GREENPLY"
$ws.Range("C15").Value = "NSE:GUJAPOLLO"
$ws.Range("C16").Value = "NSE:HEG"
$ws.Range("C17").Value = "NSE:HLEGLAS"
$ws.Range("C18").Value = "NSE:JTEKTINDIA"
$ws.Range("C19").Value = "NSE:JUBLPHARMA"
$ws.Range("C20").Value = "NSE:LPDC"
$ws.Range("C21").Value = "NSE:LXCHEM"
$ws.Range("C22").Value = "NSE:MAWANASUG"
$ws.Range("C23").Value = "NSE:MAYURUNIQ"
$ws.Range("C24").Value = "NSE:MONARCH"
$ws.Range("C25").Value = "NSE:NAVA"
$ws.Range("C26").Value = "NSE:NUCLEUS"
$ws.Range("C27").Value = "NSE:ORIENTHOT"
$ws.Range("C28").Value = "NSE:PALASHSECU"

# Add new rows 29-35
$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "NSE:PDSL"
$ws.Range("A30").Value = 28
$ws.Range("C30").Value = "NSE:PGHH"
$ws.Range("A31").Value = 29
$ws.Range("C31").Value = "NSE:PNBHOUSING"
$ws.Range("A32").Value = 30
$ws.Range("C32").Value = "NSE:PONNIERODE"
$ws.Range("A33").Value = 31
$ws.Range("C33").Value = "NSE:RELAXO"
$ws.Range("A34").Value = 32
$ws.Range("C34").Value = "NSE:RPSGVENT"
$ws.Range("A35").Value = 33
$ws.Range("C35").Value = "NSE:RRKABEL"

# Apply the bordered/centered number style (matching A2:A28) to the new index cells
$ws.Range("A2").Copy()
$ws.Range("A29:A35").PasteSpecial(-4122)
